$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

$data = @(
    @(0, 2422.766666666667, 2679, 2268, 0.02952713171641032),
    @(1, 2425, 2552, 2364, 0.02819689909617106),
    @(2, 2287.866666666667, 2411, 2171, 0.03103157679239909),
    @(3, 2155.466666666667, 2234, 2080, 0.02859067916870117),
    @(4, 2157.3, 2251, 1922, 0.03163760503133138),
    @(5, 1786.433333333333, 1949, 1614, 0.03176081975301107),
    @(6, 2599.766666666667, 2811, 2395, 0.03093178272247315),
    @(7, 2434.8, 2550, 2294, 0.03075393040974935),
    @(8, 2456.833333333333, 2761, 2082, 0.03124056657155355),
    @(9, 2268.2, 2326, 2230, 0.02914390563964844)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
